$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 634.1667
$ws.Range("I12").Value = 583.5714
$ws.Range("K12").Value = 583.5714
$ws.Range("M12").Value = -413.5714
$ws.Range("H62").Value = 2313.4
$ws.Range("I62").Value = 2149.5
$ws.Range("K62").Value = 2149.5
$ws.Range("M62").Value = -1525.5
$ws.Range("H64").Value = 2966.6667
$ws.Range("I64").Value = 2700
$ws.Range("K64").Value = 2700
$ws.Range("M64").Value = -2452
$ws.Range("H65").Value = 2313.4
$ws.Range("I65").Value = 2149.5
$ws.Range("K65").Value = 10747.5
$ws.Range("M65").Value = -7627.5
$ws.Range("H67").Value = 2966.6667
$ws.Range("I67").Value = 2700
$ws.Range("K67").Value = 2700
$ws.Range("M67").Value = -1842
$ws.Range("H74").Value = 4198.1665
$ws.Range("J74").Value = 4198.5
$ws.Range("L74").Value = 4198.5
$ws.Range("N74").Value = -6070.5
$ws.Range("H77").Value = 4198.1665
$ws.Range("J77").Value = 4198.5
$ws.Range("L77").Value = 20992.5
$ws.Range("N77").Value = -30352.5
$ws.Range("H86").Value = 2128.7144
$ws.Range("I86").Value = 2099.4
$ws.Range("K86").Value = 2099.4
$ws.Range("M86").Value = -976.4000000000001
$ws.Range("H89").Value = 2128.7144
$ws.Range("I89").Value = 2099.4
$ws.Range("K89").Value = 10497
$ws.Range("M89").Value = -4881
$ws.Range("H98").Value = 3073
$ws.Range("I98").Value = 3226.2856
$ws.Range("K98").Value = 3226.2856
$ws.Range("M98").Value = -1728.2856
$ws.Range("H112").Value = 3773.2666
$ws.Range("J112").Value = 3773.2666
$ws.Range("L112").Value = 11319.7998
$ws.Range("N112").Value = -13535.7998
$ws.Range("H122").Value = 3073
$ws.Range("I122").Value = 3226.2856
$ws.Range("K122").Value = 9678.856800000001
$ws.Range("M122").Value = -7228.856800000001
$ws.Range("H132").Value = 1134.7037
$ws.Range("I132").Value = 1134.7037
$ws.Range("K132").Value = 3404.1111
$ws.Range("M132").Value = -874.1111000000001
$ws.Range("H137").Value = 2518.3333
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 2518.3333
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 7554.999899999999
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -12654.9999
$ws.Range("H138").Value = 3005.853
$ws.Range("I138").Value = 3709.0557
$ws.Range("J138").Value = 2214.75
$ws.Range("K138").Value = 11127.1671
$ws.Range("L138").Value = 6644.25
$ws.Range("M138").Value = -5987.167099999999
$ws.Range("N138").Value = -16924.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4563.45
$ws.Range("I32").Value = 3145.8965
$ws.Range("J32").Value = 8300.637000000001
$ws.Range("K32").Value = 3145.8965
$ws.Range("L32").Value = 8300.637000000001
$ws.Range("M32").Value = -2858.8965
$ws.Range("N32").Value = -8874.637000000001
$ws.Range("H63").Value = 7374.5
$ws.Range("I63").Value = 7332.6665
$ws.Range("K63").Value = 7332.6665
$ws.Range("M63").Value = -6646.6665
$ws.Range("H66").Value = 7374.5
$ws.Range("I66").Value = 7332.6665
$ws.Range("K66").Value = 36663.3325
$ws.Range("M66").Value = -33231.3325
$ws.Range("H135").Value = 30500.125
$ws.Range("J135").Value = 30500.125
$ws.Range("L135").Value = 30500.125
$ws.Range("N135").Value = -40640.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1324.8462
$ws.Range("I20").Value = 1327.5264
$ws.Range("J20").Value = 1317.5714
$ws.Range("K20").Value = 1327.5264
$ws.Range("L20").Value = 1317.5714
$ws.Range("M20").Value = -1080.5264
$ws.Range("N20").Value = -1811.5714
$ws.Range("H81").Value = 59994.832
$ws.Range("J81").Value = 59994.832
$ws.Range("L81").Value = 59994.832
$ws.Range("N81").Value = -62116.832
$ws.Range("H82").Value = 29999.5
$ws.Range("I82").Value = 29999.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 29999.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -29616.5
$ws.Range("N82").ClearContents()
$ws.Range("H84").Value = 59994.832
$ws.Range("J84").Value = 59994.832
$ws.Range("L84").Value = 179984.496
$ws.Range("N84").Value = -190592.496
$ws.Range("H85").Value = 29999.5
$ws.Range("I85").Value = 29999.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 29999.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -28673.5
$ws.Range("N85").ClearContents()
$ws.Range("H86").Value = 155553.77
$ws.Range("I86").Value = 1833.2222
$ws.Range("J86").Value = 501425
$ws.Range("K86").Value = 1833.2222
$ws.Range("L86").Value = 501425
$ws.Range("M86").Value = -710.2221999999999
$ws.Range("N86").Value = -503671
$ws.Range("H89").Value = 155553.77
$ws.Range("I89").Value = 1833.2222
$ws.Range("J89").Value = 501425
$ws.Range("K89").Value = 9166.110999999999
$ws.Range("L89").Value = 2507125
$ws.Range("M89").Value = -3550.110999999999
$ws.Range("N89").Value = -2518357
$ws.Range("H134").Value = 11915
$ws.Range("I134").Value = 13873.214
$ws.Range("J134").Value = 5061.25
$ws.Range("K134").Value = 41619.642
$ws.Range("L134").Value = 15183.75
$ws.Range("M134").Value = -39084.642
$ws.Range("N134").Value = -20253.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H105").Value = 1959.3334
$ws.Range("I105").Value = 1951.2
$ws.Range("K105").Value = 1951.2
$ws.Range("M105").Value = -204.2
$ws.Range("H132").Value = 2463.9
$ws.Range("I132").Value = 933.4545000000001
$ws.Range("J132").Value = 4334.4443
$ws.Range("K132").Value = 2800.3635
$ws.Range("L132").Value = 13003.3329
$ws.Range("M132").Value = -270.3635000000004
$ws.Range("N132").Value = -18063.3329
$ws.Range("H134").Value = 2275.4375
$ws.Range("I134").Value = 1892.9231
$ws.Range("J134").Value = 3933
$ws.Range("K134").Value = 5678.7693
$ws.Range("L134").Value = 11799
$ws.Range("M134").Value = -3143.7693
$ws.Range("N134").Value = -16869

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 844
$ws.Range("J5").Value = 990
$ws.Range("L5").Value = 2970
$ws.Range("N5").Value = -3194
$ws.Range("H126").Value = 83336970
$ws.Range("I126").Value = 166668270
$ws.Range("K126").Value = 500004810
$ws.Range("M126").Value = -499999870
$ws.Range("H131").Value = 8487404
$ws.Range("J131").Value = 14748.157
$ws.Range("L131").Value = 44244.471
$ws.Range("N131").Value = -54324.471
$ws.Range("H135").Value = 844
$ws.Range("J135").Value = 990
$ws.Range("L135").Value = 8910
$ws.Range("N135").Value = -13980

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3853
$ws.Range("I70").Value = 4059.5557
$ws.Range("K70").Value = 4059.5557
$ws.Range("M70").Value = -3789.5557
$ws.Range("H73").Value = 3853
$ws.Range("I73").Value = 4059.5557
$ws.Range("K73").Value = 4059.5557
$ws.Range("M73").Value = -3123.5557
$ws.Range("H102").Value = 2095.818
$ws.Range("I102").Value = 2301.818
$ws.Range("K102").Value = 2301.818
$ws.Range("M102").Value = -679.8180000000002
$ws.Range("H122").Value = 2484.2856
$ws.Range("I122").Value = 2278
$ws.Range("K122").Value = 6834
$ws.Range("M122").Value = -4384

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8114.75
$ws.Range("J16").Value = 500
$ws.Range("L16").Value = 500
$ws.Range("N16").Value = -840
$ws.Range("I74").Value = 50000
$ws.Range("K74").Value = 50000
$ws.Range("M74").Value = -49002
$ws.Range("I77").Value = 50000
$ws.Range("K77").Value = 150000
$ws.Range("M77").Value = -145008
$ws.Range("H136").Value = 5052.125
$ws.Range("I136").Value = 3500.5
$ws.Range("J136").Value = 5569.3335
$ws.Range("K136").Value = 10501.5
$ws.Range("L136").Value = 16708.0005
$ws.Range("M136").Value = -7951.5
$ws.Range("N136").Value = -21808.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 27999.5
$ws.Range("J70").Value = 27999.5
$ws.Range("L70").Value = 27999.5
$ws.Range("N70").Value = -28629.5
$ws.Range("H73").Value = 27999.5
$ws.Range("J73").Value = 27999.5
$ws.Range("L73").Value = 27999.5
$ws.Range("N73").Value = -30183.5
$ws.Range("H100").Value = 446.66666
$ws.Range("I100").Value = 336
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 672
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -131
$ws.Range("N100").Value = -3082
$ws.Range("H132").Value = 1710.4
$ws.Range("I132").Value = 954
$ws.Range("K132").Value = 2862
$ws.Range("M132").Value = -332
